$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues = -4163
$xlPasteValues = -4163

$ws.Range("D2").Value = "29.834.41"
$ws.Range("E2").Value = "  -1.25%  "

$ws.Range("D3").Value = "1.892.14"
$ws.Range("E3").Value = "  -1.08%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Formula = "=""0.7818"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial($xlPasteValues)
$ws.Range("E5").Value = "  -4.86%  "

$ws.Range("D6").Formula = "=""243.69"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial($xlPasteValues)
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Formula = "=""0.3139"""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial($xlPasteValues)
$ws.Range("E8").Value = "  -3.66%  "

$ws.Range("D9").Formula = "=""25.26"""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial($xlPasteValues)
$ws.Range("E9").Value = "  -5.96%  "

$ws.Range("D10").Formula = "=""0.07171"""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial($xlPasteValues)
$ws.Range("E10").Value = "  +1.51%  "

$ws.Range("D11").Formula = "=""0.08088"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial($xlPasteValues)
$ws.Range("E11").Value = "  -0.16%  "

$ws.Range("D12").Formula = "=""0.7629"""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial($xlPasteValues)
$ws.Range("E12").Value = "  -1.45%  "

$ws.Range("D13").Formula = "=""5.483"""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial($xlPasteValues)
$ws.Range("E13").Value = "  +3.52%  "

$ws.Range("D14").Value = "1.848.19"
$ws.Range("E14").Value = "  -3.37%  "

$ws.Range("D15").Formula = "=""92.25"""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial($xlPasteValues)
$ws.Range("E15").Value = "  -1.29%  "

$ws.Range("D16").Formula = "=""6.134"""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial($xlPasteValues)
$ws.Range("E16").Value = "  +3.49%  "

$ws.Range("D17").Value = "29.823.07"
$ws.Range("E17").Value = "  -1.26%  "

$ws.Range("D18").Formula = "=""13.95"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial($xlPasteValues)
$ws.Range("E18").Value = "  -2.10%  "

$ws.Range("D19").Formula = "=""243.02"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial($xlPasteValues)
$ws.Range("E19").Value = "  -1.50%  "

$ws.Range("D20").Formula = "=""0.000007774"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial($xlPasteValues)
$ws.Range("E20").Value = "  -0.38%  "

$ws.Range("D21").Formula = "=""1.001"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial($xlPasteValues)

$ws.Range("D22").Value = "2.148.34"
$ws.Range("E22").Value = "  -0.78%  "

$ws.Range("D23").Formula = "=""8.105"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial($xlPasteValues)
$ws.Range("E23").Value = "  +14.03%  "

$ws.Range("D24").Formula = "=""1.002"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial($xlPasteValues)

$ws.Range("D25").Formula = "=""0.1633"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial($xlPasteValues)
$ws.Range("E25").Value = "  -2.64%  "

$ws.Range("D26").Formula = "=""9.407"""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial($xlPasteValues)
$ws.Range("E26").Value = "  +0.70%  "

$ws.Range("D27").Formula = "=""162.71"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial($xlPasteValues)
$ws.Range("E27").Value = "  -2.74%  "

$ws.Range("D28").Formula = "=""18.71"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial($xlPasteValues)
$ws.Range("E28").Value = "  -1.56%  "

$ws.Range("D29").Formula = "=""2.046"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial($xlPasteValues)
$ws.Range("E29").Value = "  -3.23%  "

$ws.Range("D30").Formula = "=""1.411"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial($xlPasteValues)
$ws.Range("E30").Value = "  +2.66%  "

$ws.Range("D31").Formula = "=""1.547"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial($xlPasteValues)
$ws.Range("E31").Value = "  +1.17%  "

$ws.Range("D32").Formula = "=""4.485"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial($xlPasteValues)
$ws.Range("E32").Value = "  +3.95%  "

$ws.Range("D33").Formula = "=""4.112"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial($xlPasteValues)
$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("D34").Formula = "=""0.05549"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial($xlPasteValues)
$ws.Range("E34").Value = "  -5.40%  "

$ws.Range("D35").Formula = "=""1.265"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial($xlPasteValues)
$ws.Range("E35").Value = "  -0.94%  "

$ws.Range("D36").Formula = "=""0.7427"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial($xlPasteValues)
$ws.Range("E36").Value = "  +0.71%  "

$ws.Range("D37").Formula = "=""0.9965"""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial($xlPasteValues)
$ws.Range("E37").Value = "  -0.32%  "

$ws.Range("D38").Formula = "=""2.617"""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial($xlPasteValues)
$ws.Range("E38").Value = "  -2.85%  "

$ws.Range("D39").Formula = "=""0.01918"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial($xlPasteValues)
$ws.Range("E39").Value = "  -0.42%  "

$ws.Range("D40").Formula = "=""2.785"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial($xlPasteValues)
$ws.Range("E40").Value = "  -0.55%  "

$ws.Range("D41").Value = "1.146.57"
$ws.Range("E41").Value = "  +13.56%  "

$ws.Range("D42").Formula = "=""73.59"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial($xlPasteValues)
$ws.Range("E42").Value = "  +0.14%  "

$ws.Range("D43").Formula = "=""0.4409"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial($xlPasteValues)
$ws.Range("E43").Value = "  -1.38%  "

$ws.Range("E44").Value = "  -1.98%  "

$ws.Range("E45").Value = "  -0.35%  "

$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("D47").Formula = "=""103.66"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial($xlPasteValues)
$ws.Range("E47").Value = "  +0.75%  "

$ws.Range("D48").Formula = "=""1.874"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial($xlPasteValues)
$ws.Range("E48").Value = "  -2.27%  "

$ws.Range("D49").Formula = "=""9.957"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial($xlPasteValues)
$ws.Range("E49").Value = "  +0.76%  "

$ws.Range("D50").Formula = "=""7.442"""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial($xlPasteValues)
$ws.Range("E50").Value = "  -2.13%  "

$ws.Range("D51").Formula = "=""2.996"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial($xlPasteValues)
$ws.Range("E51").Value = "  +9.76%  "
